$d = $word.ActiveDocument

# The TODO list currently has three separate bullet paragraphs:
#   "Create pink ghost behavior."
#   "Create Blue ghost behavior."
#   "Create Yellow ghost behavior."  (this last one also carries the _GoBack bookmark)
#
# We need to collapse them into a single bullet paragraph with new text,
# while preserving the bookmark that lives at the end of the third paragraph.
#
# Merging paragraphs by deleting just the paragraph-mark character (rather
# than deleting a whole paragraph Range, and rather than replacing across
# "^p" with Find) keeps the trailing bookmark intact.

$pinkPara = $null
$bluePara = $null
$yellowPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -eq "Create pink ghost behavior.`r") { $pinkPara = $i }
    elseif ($txt -eq "Create Blue ghost behavior.`r") { $bluePara = $i }
    elseif ($txt -eq "Create Yellow ghost behavior.`r") { $yellowPara = $i }
}

# Delete the paragraph mark ending the "Blue" paragraph first (higher index),
# merging "Blue" and "Yellow" (with its bookmark) into one paragraph.
$pBlue = $d.Paragraphs.Item($bluePara)
$markBlue = $d.Range($pBlue.Range.End - 1, $pBlue.Range.End)
$markBlue.Delete()

# Delete the paragraph mark ending the "pink" paragraph, merging it with
# the combined "Blue"+"Yellow" paragraph.
$pPink = $d.Paragraphs.Item($pinkPara)
$markPink = $d.Range($pPink.Range.End - 1, $pPink.Range.End)
$markPink.Delete()

# Now replace the concatenated run text with the new single sentence. The
# bookmark (now trailing in this same paragraph) is untouched by this Find.
$pPink = $d.Paragraphs.Item($pinkPara)
$pPink.Range.Find.Execute(
    "Create pink ghost behavior.Create Blue ghost behavior.Create Yellow ghost behavior.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Fix eyes not facing the right way on animation for every ghost except for red ghost.",
    2) | Out-Null
